$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19 (Leve Item ID 7015)
$ws.Range("H19").Value = 692.93335
$ws.Range("I19").Value = 750
$ws.Range("J19").Value = 684.1539
$ws.Range("K19").Value = 750
$ws.Range("L19").Value = 684.1539
$ws.Range("M19").Value = -575
$ws.Range("N19").Value = -1034.1539

# Row 34 (Leve Item ID 2160)
$ws.Range("H34").Value = 1886.3334
$ws.Range("I34").Value = 1886.3334
$ws.Range("K34").Value = 1886.3334
$ws.Range("M34").Value = -1683.3334

# Row 36 (Leve Item ID 2160)
$ws.Range("H36").Value = 1886.3334
$ws.Range("I36").Value = 1886.3334
$ws.Range("K36").Value = 1886.3334
$ws.Range("M36").Value = -1171.3334

# Row 51 (Leve Item ID 5486)
$ws.Range("H51").Value = 4770.2
$ws.Range("I51").Value = 3933.3333
$ws.Range("J51").Value = 5128.857
$ws.Range("K51").Value = 3933.3333
$ws.Range("L51").Value = 5128.857
$ws.Range("M51").Value = -3449.3333
$ws.Range("N51").Value = -6096.857

# Row 103 (Leve Item ID 19909)
$ws.Range("H103").Value = 489.25
$ws.Range("I103").Value = 344.85715
$ws.Range("K103").Value = 1034.57145
$ws.Range("M103").Value = -448.5714499999999

$ws = $wb.Worksheets.Item("ARM")
# Row 43 (Leve Item ID 21715)
$ws.Range("H43").Value = 7981.6665
$ws.Range("J43").Value = 7981.6665
$ws.Range("L43").Value = 7981.6665
$ws.Range("N43").Value = -8607.666499999999

# Row 45 (Leve Item ID 27714)
$ws.Range("H45").Value = 2883.2856
$ws.Range("I45").Value = 3114.1667
$ws.Range("J45").Value = 1498
$ws.Range("K45").Value = 3114.1667
$ws.Range("L45").Value = 1498
$ws.Range("M45").Value = -2737.1667
$ws.Range("N45").Value = -2252

# Row 88 (Leve Item ID 12530)
$ws.Range("H88").Value = 6411.1177
$ws.Range("J88").Value = 8137.1816
$ws.Range("L88").Value = 8137.1816
$ws.Range("N88").Value = -8949.1816

# Row 91 (Leve Item ID 12530)
$ws.Range("H91").Value = 6411.1177
$ws.Range("J91").Value = 8137.1816
$ws.Range("L91").Value = 8137.1816
$ws.Range("N91").Value = -10945.1816

# Row 110 (Leve Item ID 27708)
$ws.Range("H110").Value = 2199.2856
$ws.Range("I110").Value = 3122.75
$ws.Range("J110").Value = 968
$ws.Range("K110").Value = 3122.75
$ws.Range("L110").Value = 968
$ws.Range("M110").Value = -1077.75
$ws.Range("N110").Value = -5058

$ws = $wb.Worksheets.Item("BSM")
# Row 86 (Leve Item ID 12526)
$ws.Range("H86").Value = 13430.6
$ws.Range("I86").Value = 23950.2
$ws.Range("J86").Value = 2911
$ws.Range("K86").Value = 23950.2
$ws.Range("L86").Value = 2911
$ws.Range("M86").Value = -22827.2
$ws.Range("N86").Value = -5157

# Row 89 (Leve Item ID 12526)
$ws.Range("H89").Value = 13430.6
$ws.Range("I89").Value = 23950.2
$ws.Range("J89").Value = 2911
$ws.Range("K89").Value = 119751
$ws.Range("L89").Value = 14555
$ws.Range("M89").Value = -114135
$ws.Range("N89").Value = -25787

# Row 105 (Leve Item ID 19947)
$ws.Range("H105").Value = 22729756
$ws.Range("I105").Value = 35716332
$ws.Range("K105").Value = 35716332
$ws.Range("M105").Value = -35714585

# Row 107 (Leve Item ID 27706)
$ws.Range("H107").Value = 3501.625
$ws.Range("I107").Value = 4500.5
$ws.Range("J107").Value = 3168.6667
$ws.Range("K107").Value = 4500.5
$ws.Range("L107").Value = 3168.6667
$ws.Range("M107").Value = -2580.5
$ws.Range("N107").Value = -7008.6667

$ws = $wb.Worksheets.Item("CRP")
# Row 5 (Leve Item ID 1893)
$ws.Range("H5").Value = 494.54544
$ws.Range("I5").Value = 81.333336
$ws.Range("J5").Value = 649.5
$ws.Range("K5").Value = 81.333336
$ws.Range("L5").Value = 649.5
$ws.Range("M5").Value = 30.666664
$ws.Range("N5").Value = -873.5

# Row 16 (Leve Item ID 27691)
$ws.Range("H16").Value = 111113450
$ws.Range("I16").Value = 2000
$ws.Range("J16").Value = 166669170
$ws.Range("K16").Value = 2000
$ws.Range("L16").Value = 166669170
$ws.Range("M16").Value = -1713
$ws.Range("N16").Value = -166669744

# Row 25 (Leve Item ID 1895)
$ws.Range("H25").Value = 1000
$ws.Range("I25").Value = 1000
$ws.Range("K25").Value = 1000
$ws.Range("M25").Value = -826

# Row 35 (Leve Item ID 1627)
$ws.Range("H35").Value = 520
$ws.Range("I35").Value = 520
$ws.Range("K35").Value = 520
$ws.Range("M35").Value = -226

# Row 41 (Leve Item ID 1917)
$ws.Range("H41").Value = 8660
$ws.Range("I41").Value = 1650
$ws.Range("J41").Value = 13333.333
$ws.Range("K41").Value = 1650
$ws.Range("L41").Value = 13333.333
$ws.Range("M41").Value = -1222
$ws.Range("N41").Value = -14189.333

# Row 59 (Leve Item ID 1942)
$ws.Range("H59").Value = 30063.5
$ws.Range("J59").Value = 30063.5
$ws.Range("L59").Value = 30063.5
$ws.Range("N59").Value = -32353.5

# Row 60 (Leve Item ID 1937)
$ws.Range("H60").Value = 11103
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 11103
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 11103
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = -12125

# Row 74 (Leve Item ID 10636)
$ws.Range("H74").Value = 25314
$ws.Range("J74").Value = 25314
$ws.Range("L74").Value = 25314
$ws.Range("N74").Value = -27062

# Row 77 (Leve Item ID 10636)
$ws.Range("H77").Value = 25314
$ws.Range("J77").Value = 25314
$ws.Range("L77").Value = 75942
$ws.Range("N77").Value = -84678

# Row 105 (Leve Item ID 19928)
$ws.Range("H105").Value = 1356.3334
$ws.Range("I105").Value = 784.5
$ws.Range("J105").Value = 2500
$ws.Range("K105").Value = 784.5
$ws.Range("L105").Value = 2500
$ws.Range("M105").Value = 962.5
$ws.Range("N105").Value = -5994

# Row 107 (Leve Item ID 27689)
$ws.Range("H107").Value = 1312.5
$ws.Range("I107").Value = 850
$ws.Range("J107").Value = 1775
$ws.Range("K107").Value = 850
$ws.Range("L107").Value = 1775
$ws.Range("M107").Value = 1070
$ws.Range("N107").Value = -5615

# Row 113 (Leve Item ID 27691)
$ws.Range("H113").Value = 111113450
$ws.Range("I113").Value = 2000
$ws.Range("J113").Value = 166669170
$ws.Range("K113").Value = 2000
$ws.Range("L113").Value = 166669170
$ws.Range("M113").Value = 170
$ws.Range("N113").Value = -166673510

$ws = $wb.Worksheets.Item("CUL")
# Row 80 (Leve Item ID 12890)
$ws.Range("H80").Value = 4200
$ws.Range("I80").Value = 3300
$ws.Range("K80").Value = 9900
$ws.Range("M80").Value = -8964

# Row 83 (Leve Item ID 12890)
$ws.Range("H83").Value = 4200
$ws.Range("I83").Value = 3300
$ws.Range("K83").Value = 29700
$ws.Range("M83").Value = -25020

# Row 92 (Leve Item ID 19841)
$ws.Range("H92").Value = 1350
$ws.Range("I92").Value = 300
$ws.Range("J92").Value = 1560
$ws.Range("K92").Value = 900
$ws.Range("L92").Value = 4680
$ws.Range("M92").Value = 348
$ws.Range("N92").Value = -7176

# Row 113 (Leve Item ID 27843)
$ws.Range("H113").Value = 562.8913
$ws.Range("I113").Value = 525.3333
$ws.Range("J113").Value = 568.525
$ws.Range("K113").Value = 1575.9999
$ws.Range("L113").Value = 1705.575
$ws.Range("M113").Value = 594.0001
$ws.Range("N113").Value = -6045.575

$ws = $wb.Worksheets.Item("GSM")
# Row 80 (Leve Item ID 12521)
$ws.Range("H80").Value = 2852.3333
$ws.Range("I80").Value = 2103.6667
$ws.Range("J80").Value = 3601
$ws.Range("K80").Value = 2103.6667
$ws.Range("L80").Value = 3601
$ws.Range("M80").Value = -1105.6667
$ws.Range("N80").Value = -5597

# Row 83 (Leve Item ID 12521)
$ws.Range("H83").Value = 2852.3333
$ws.Range("I83").Value = 2103.6667
$ws.Range("J83").Value = 3601
$ws.Range("K83").Value = 10518.3335
$ws.Range("L83").Value = 18005
$ws.Range("M83").Value = -5526.333500000001
$ws.Range("N83").Value = -27989

$ws = $wb.Worksheets.Item("LTW")
# Row 18 (Leve Item ID 3772)
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()

# Row 40 (Leve Item ID 36248)
$ws.Range("H40").Value = 2332.2144
$ws.Range("I40").Value = 2220.9167
$ws.Range("J40").Value = 3000
$ws.Range("K40").Value = 2220.9167
$ws.Range("L40").Value = 3000
$ws.Range("M40").Value = -2084.9167
$ws.Range("N40").Value = -3272

$ws = $wb.Worksheets.Item("WVR")
# Row 102 (Leve Item ID 19642)
$ws.Range("H102").Value = 30000
$ws.Range("J102").Value = 30000
$ws.Range("L102").Value = 30000
$ws.Range("N102").Value = -36490

# Row 107 (Leve Item ID 27746)
$ws.Range("H107").Value = 451.75
$ws.Range("I107").Value = 402
$ws.Range("J107").Value = 501.5
$ws.Range("K107").Value = 1206
$ws.Range("L107").Value = 1504.5
$ws.Range("M107").Value = 714
$ws.Range("N107").Value = -5344.5
